$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data set (rows 2..22), each entry: ResidueCombination, Count, PermeatedResidues, Frames
$data = @(
    @("98, 130, 748, 1073", 3, "130, 130, 130", "5178, 5582, 6488"),
    @("98, 455, 780, 1105", 1, "780", "5552"),
    @("130, 780, 780", 1, "780", "3171"),
    @("98, 130, 748, 780", 1, "130", "3631"),
    @("98, 780, 1105", 1, "780", "3666"),
    @("130, 455, 748, 1073", 1, "130", "5131"),
    @("98, 130, 1073, SF", 1, "130", "4415"),
    @("423, 1073, 1105, SF, SF", 1, "1105", "4994"),
    @("130, 455, 748, 780", 1, "130", "5269"),
    @("130, 748, 1073, 1105", 1, "1105", "5399"),
    @("423, 748, 780, 1073", 1, "780", "5677"),
    @("423, 748, 1073, 1105", 2, "1105, 1105", "5331, 5433"),
    @("130, 423, 748, 1073", 1, "130", "5886"),
    @("98, 130, 455, 748, 1073", 1, "130", "6016"),
    @("130, 423, 748, 1073, SF", 1, "130", "6202"),
    @("130, 130, 423", 1, "130", "6426"),
    @("130, 455, 780, 1105", 1, "130", "6561"),
    @("98, 780, 780, 1105", 1, "780", "6359"),
    @("98, 130, 130, 455, 780", 1, "130", "6727"),
    @("130, 130, 423, 1073", 1, "130", "6670"),
    @("98, 98, 455, 455", 1, "455", "6748")
)

$rowIndex = 2
foreach ($entry in $data) {
    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.Value = "'" + $entry[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($rowIndex, 2).Value = $entry[1]

    $cellC = $ws.Cells.Item($rowIndex, 3)
    $cellC.Value = "'" + $entry[2]
    $cellC.Style = "Normal"

    $cellD = $ws.Cells.Item($rowIndex, 4)
    $cellD.Value = "'" + $entry[3]
    $cellD.Style = "Normal"

    $rowIndex++
}
